$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 20
# from serial date 45192 to 45202 (10 days later).
for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}
